$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 1, 0.03888888888888889, 0.06458333333333334, "E20", 3),
    @(3, 2, 0.01041666666666667, 0.03402777777777777, "C22", 1),
    @(4, 3, 0.03611111111111111, 0.06319444444444444, "C18", 1),
    @(5, 4, 0.02777777777777778, 0.05069444444444444, "D31", 1),
    @(6, 5, 0.03263888888888889, 0.05833333333333333, "F57", 1),
    @(7, 6, 0.01388888888888889, 0.03680555555555556, "D36", 2),
    @(8, 7, 0.02708333333333333, 0.05069444444444444, "F56", 1),
    @(9, 8, 0.03611111111111111, 0.05833333333333333, "A18", 1),
    @(10, 9, 0.01875, 0.04166666666666666, "D33", 1),
    @(11, 10, 0.002083333333333333, 0.02638888888888889, "A14", 3),
    @(12, 11, 0.00625, 0.03194444444444444, "D44", 3),
    @(13, 12, 0.06597222222222222, 0.09027777777777778, "A5", 2),
    @(14, 13, 0.07847222222222222, 0.1048611111111111, "A19", 1),
    @(15, 14, 0.04236111111111111, 0.06527777777777778, "F54", 3),
    @(16, 15, 0.06111111111111111, 0.08541666666666667, "F31", 2),
    @(17, 16, 0.05625, 0.07708333333333334, "C26", 1),
    @(18, 17, 0.05416666666666667, 0.0798611111111111, "F60", 3),
    @(19, 18, 0.07916666666666666, 0.1034722222222222, "D44", 1),
    @(20, 19, 0.05416666666666667, 0.075, "A8", 1),
    @(21, 20, 0.07708333333333334, 0.1006944444444444, "A9", 1),
    @(22, 21, 0.04375, 0.06527777777777778, "A21", 3),
    @(23, 22, 0.075, 0.1020833333333333, "B7", 2),
    @(24, 23, 0.1, 0.1222222222222222, "E1", 1),
    @(25, 24, 0.1027777777777778, 0.1236111111111111, "F60", 3),
    @(26, 25, 0.08402777777777778, 0.1097222222222222, "E7", 2),
    @(27, 26, 0.1194444444444445, 0.1402777777777778, "E27", 1),
    @(28, 27, 0.1236111111111111, 0.1472222222222222, "C15", 3),
    @(29, 28, 0.10625, 0.13125, "B7", 1),
    @(30, 29, 0.08333333333333333, 0.10625, "F34", 1),
    @(31, 30, 0.1097222222222222, 0.13125, "B1", 2),
    @(32, 31, 0.1097222222222222, 0.1326388888888889, "B8", 2),
    @(33, 32, 0.1152777777777778, 0.1395833333333333, "D40", 3),
    @(34, 33, 0.1159722222222222, 0.1409722222222222, "F52", 2),
    @(35, 34, 0.1152777777777778, 0.1361111111111111, "D42", 1),
    @(36, 35, 0.11875, 0.1430555555555555, "F55", 1),
    @(37, 36, 0.08819444444444445, 0.1138888888888889, "D34", 1),
    @(38, 37, 0.1111111111111111, 0.1381944444444445, "D49", 1),
    @(39, 38, 0.09652777777777778, 0.1201388888888889, "C17", 2),
    @(40, 39, 0.08680555555555555, 0.1118055555555556, "D37", 1),
    @(41, 40, 0.08611111111111111, 0.1131944444444444, "A7", 2),
    @(42, 41, 0.1395833333333333, 0.1625, "F35", 2),
    @(43, 42, 0.1430555555555555, 0.1645833333333333, "D48", 1),
    @(44, 43, 0.1534722222222222, 0.1756944444444444, "A20", 1),
    @(45, 44, 0.13125, 0.1541666666666667, "C1", 2),
    @(46, 45, 0.1256944444444444, 0.1479166666666667, "C21", 2),
    @(47, 46, 0.1375, 0.1638888888888889, "B5", 1),
    @(48, 47, 0.1375, 0.1645833333333333, "F51", 3),
    @(49, 48, 0.1284722222222222, 0.15, "A7", 3),
    @(50, 49, 0.1645833333333333, 0.1881944444444444, "A18", 3),
    @(51, 50, 0.1354166666666667, 0.1569444444444444, "A9", 1),
    @(52, 51, 0.1263888888888889, 0.1479166666666667, "A2", 3),
    @(53, 52, 0.1486111111111111, 0.1708333333333333, "D37", 1),
    @(54, 53, 0.1298611111111111, 0.1548611111111111, "A8", 3),
    @(55, 54, 0.1479166666666667, 0.1743055555555555, "E6", 1),
    @(56, 55, 0.1604166666666667, 0.1868055555555556, "D36", 2),
    @(57, 56, 0.1979166666666667, 0.2229166666666667, "A13", 3),
    @(58, 57, 0.19375, 0.2159722222222222, "E12", 1),
    @(59, 58, 0.1923611111111111, 0.2173611111111111, "E28", 1),
    @(60, 59, 0.1819444444444444, 0.2076388888888889, "E22", 1),
    @(61, 60, 0.1798611111111111, 0.2027777777777778, "E8", 2),
    @(62, 61, 0.1951388888888889, 0.2222222222222222, "C16", 3),
    @(63, 62, 0.1666666666666667, 0.1930555555555556, "A17", 3),
    @(64, 63, 0.1986111111111111, 0.2194444444444444, "D35", 3),
    @(65, 64, 0.1833333333333333, 0.2097222222222222, "F54", 1),
    @(66, 65, 0.2145833333333333, 0.2388888888888889, "A17", 3),
    @(67, 66, 0.2159722222222222, 0.2430555555555556, "A21", 3),
    @(68, 67, 0.2388888888888889, 0.2625, "D48", 3),
    @(69, 68, 0.2465277777777778, 0.2673611111111111, "D46", 3),
    @(70, 69, 0.2083333333333333, 0.2326388888888889, "E27", 3),
    @(71, 70, 0.2458333333333333, 0.2666666666666667, "A21", 2),
    @(72, 71, 0.21875, 0.2395833333333333, "E26", 2),
    @(73, 72, 0.21875, 0.2402777777777778, "B7", 1),
    @(74, 73, 0.2090277777777778, 0.2340277777777778, "A4", 2),
    @(75, 74, 0.2222222222222222, 0.2472222222222222, "C25", 2),
    @(76, 75, 0.2618055555555556, 0.2881944444444444, "A18", 1),
    @(77, 76, 0.2909722222222222, 0.3173611111111111, "F53", 2),
    @(78, 77, 0.2694444444444444, 0.2902777777777778, "A13", 1),
    @(79, 78, 0.2798611111111111, 0.30625, "D36", 2),
    @(80, 79, 0.2645833333333333, 0.2916666666666667, "E8", 3),
    @(81, 80, 0.2770833333333333, 0.2993055555555555, "D49", 2),
    @(82, 81, 0.2625, 0.2868055555555555, "A15", 3),
    @(83, 82, 0.3159722222222222, 0.3381944444444445, "F41", 3),
    @(84, 83, 0.2944444444444445, 0.3194444444444444, "E21", 2),
    @(85, 84, 0.3145833333333333, 0.3416666666666667, "D35", 1),
    @(86, 85, 0.3020833333333333, 0.3270833333333333, "C18", 2),
    @(87, 86, 0.2993055555555555, 0.3208333333333334, "C26", 3),
    @(88, 87, 0.3222222222222222, 0.3493055555555555, "A13", 2),
    @(89, 88, 0.3083333333333333, 0.3354166666666666, "E3", 3),
    @(90, 89, 0.3284722222222222, 0.3513888888888889, "C15", 3),
    @(91, 90, 0.3430555555555556, 0.3701388888888889, "E3", 3),
    @(92, 91, 0.3708333333333333, 0.3979166666666666, "E11", 3),
    @(93, 92, 0.3604166666666667, 0.3875, "F32", 3),
    @(94, 93, 0.3631944444444444, 0.3854166666666667, "E10", 1),
    @(95, 94, 0.3430555555555556, 0.3659722222222222, "C2", 2),
    @(96, 95, 0.3361111111111111, 0.3583333333333333, "E7", 2),
    @(97, 96, 0.3444444444444444, 0.3708333333333333, "F53", 2),
    @(98, 97, 0.4131944444444444, 0.4381944444444444, "F40", 2),
    @(99, 98, 0.3875, 0.4097222222222222, "A18", 1),
    @(100, 99, 0.4020833333333333, 0.4256944444444444, "C16", 3),
    @(101, 100, 0.3756944444444444, 0.4013888888888889, "D42", 3)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
    $ws.Cells.Item($r, 5).Value = $entry[4]
    $ws.Cells.Item($r, 6).Value = $entry[5]
}
